{"js": "// Replace the answer text in each table cell with the newly generated\n// two-digit-divided-by-one-digit answer, matching the upstream diff.\n// Each (find -> replace) pair is unique across the document, so a\n// simple exact, case-sensitive search-and-replace per pair is safe.\nconst replacements = [\n  [\"46\u00f77=6, 4\", \"67\u00f79=7, 4\"],\n  [\"35\u00f79=3, 8\", \"50\u00f76=8, 2\"],\n  [\"84\u00f78=10, 4\", \"93\u00f79=10, 3\"],\n  [\"59\u00f72=29, 1\", \"99\u00f74=24, 3\"],\n  [\"76\u00f77=10, 6\", \"22\u00f74=5, 2\"],\n  [\"14\u00f72=7, 0\", \"84\u00f74=21, 0\"],\n  [\"45\u00f76=7, 3\", \"11\u00f73=3, 2\"],\n  [\"75\u00f76=12, 3\", \"43\u00f76=7, 1\"],\n  [\"83\u00f75=16, 3\", \"84\u00f78=10, 4\"],\n  [\"82\u00f75=16, 2\", \"77\u00f75=15, 2\"],\n  [\"31\u00f76=5, 1\", \"53\u00f79=5, 8\"],\n  [\"34\u00f73=11, 1\", \"52\u00f77=7, 3\"],\n  [\"49\u00f77=7, 0\", \"20\u00f72=10, 0\"],\n  [\"41\u00f78=5, 1\", \"84\u00f72=42, 0\"],\n  [\"80\u00f79=8, 8\", \"78\u00f76=13, 0\"],\n  [\"97\u00f73=32, 1\", \"45\u00f79=5, 0\"],\n  [\"38\u00f72=19, 0\", \"89\u00f74=22, 1\"],\n  [\"80\u00f74=20, 0\", \"26\u00f78=3, 2\"],\n  [\"73\u00f72=36, 1\", \"40\u00f77=5, 5\"],\n  [\"41\u00f77=5, 6\", \"96\u00f72=48, 0\"],\n  [\"89\u00f78=11, 1\", \"75\u00f72=37, 1\"],\n  [\"50\u00f72=25, 0\", \"53\u00f74=13, 1\"],\n  [\"45\u00f72=22, 1\", \"43\u00f76=7, 1\"],\n  [\"79\u00f78=9, 7\", \"36\u00f73=12, 0\"],\n];\n\nfor (const [findText, replaceText] of replacements) {\n  const results = context.document.body.search(findText, {\n    matchCase: true,\n    matchWholeWord: true,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${findText}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(replaceText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the answer text in each table cell with the newly generated\n# two-digit-divided-by-one-digit answer, matching the upstream diff.\n# Each (find -> replace) pair is unique across the document, so a\n# simple exact, case-sensitive Find/Replace per pair is safe.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"46\u00f77=6, 4\", \"67\u00f79=7, 4\"),\n    @(\"35\u00f79=3, 8\", \"50\u00f76=8, 2\"),\n    @(\"84\u00f78=10, 4\", \"93\u00f79=10, 3\"),\n    @(\"59\u00f72=29, 1\", \"99\u00f74=24, 3\"),\n    @(\"76\u00f77=10, 6\", \"22\u00f74=5, 2\"),\n    @(\"14\u00f72=7, 0\", \"84\u00f74=21, 0\"),\n    @(\"45\u00f76=7, 3\", \"11\u00f73=3, 2\"),\n    @(\"75\u00f76=12, 3\", \"43\u00f76=7, 1\"),\n    @(\"83\u00f75=16, 3\", \"84\u00f78=10, 4\"),\n    @(\"82\u00f75=16, 2\", \"77\u00f75=15, 2\"),\n    @(\"31\u00f76=5, 1\", \"53\u00f79=5, 8\"),\n    @(\"34\u00f73=11, 1\", \"52\u00f77=7, 3\"),\n    @(\"49\u00f77=7, 0\", \"20\u00f72=10, 0\"),\n    @(\"41\u00f78=5, 1\", \"84\u00f72=42, 0\"),\n    @(\"80\u00f79=8, 8\", \"78\u00f76=13, 0\"),\n    @(\"97\u00f73=32, 1\", \"45\u00f79=5, 0\"),\n    @(\"38\u00f72=19, 0\", \"89\u00f74=22, 1\"),\n    @(\"80\u00f74=20, 0\", \"26\u00f78=3, 2\"),\n    @(\"73\u00f72=36, 1\", \"40\u00f77=5, 5\"),\n    @(\"41\u00f77=5, 6\", \"96\u00f72=48, 0\"),\n    @(\"89\u00f78=11, 1\", \"75\u00f72=37, 1\"),\n    @(\"50\u00f72=25, 0\", \"53\u00f74=13, 1\"),\n    @(\"45\u00f72=22, 1\", \"43\u00f76=7, 1\"),\n    @(\"79\u00f78=9, 7\", \"36\u00f73=12, 0\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Execute($findText, $false, $true, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n\n$d.Saved = $false\n"}
